# V. 69 — add "La huella del mal" and "Alimañas" to the "Películas" table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")
$lo = $ws.ListObjects.Item("Tabla24")

# --- Add "La huella del mal" -------------------------------------------------
$row1 = $lo.ListRows.Add()
$r1 = $row1.Range.Row
$ws.Cells.Item($r1, 2).Value = "La huella del mal"
$ws.Cells.Item($r1, 4).Value = 6
$ws.Cells.Item($r1, 5).Value = 6
$ws.Cells.Item($r1, 6).Value = 6
$ws.Cells.Item($r1, 7).Value = 6
$ws.Cells.Item($r1, 8).Value = 4.8
$ws.Cells.Item($r1, 9).Value = 4.3
$ws.Cells.Item($r1, 3).Formula = "=AVERAGE(D$r1,E$r1,E$r1,F$r1,G$r1,H$r1,H$r1,I$r1)"

# --- Add "Alimañas" ----------------------------------------------------------
$row2 = $lo.ListRows.Add()
$r2 = $row2.Range.Row
$ws.Cells.Item($r2, 2).Value = "Alimañas"
$ws.Cells.Item($r2, 4).Value = 4
$ws.Cells.Item($r2, 5).Value = 2
$ws.Cells.Item($r2, 6).Value = 6
$ws.Cells.Item($r2, 7).Value = 6
$ws.Cells.Item($r2, 8).Value = 4.9
$ws.Cells.Item($r2, 9).Value = 4.4
$ws.Cells.Item($r2, 3).Formula = "=AVERAGE(D$r2,E$r2,E$r2,F$r2,G$r2,H$r2,H$r2,I$r2)"

# Copy the normal row formatting (number formats / alignment) down onto the
# two freshly-appended rows so they match the rest of the table body.
$lastExisting = $r1 - 1
$ws.Range("B" + $lastExisting + ":I" + $lastExisting).Copy()
$ws.Range("B" + $r1 + ":I" + $r1).PasteSpecial(-4122)
$ws.Range("B" + $r2 + ":I" + $r2).PasteSpecial(-4122)

# Highlight the two new entries the same way the previous "latest addition"
# (Rebel Ridge, row 43) was highlighted, then clear that old highlight since
# it is no longer the newest entry.
$ws.Range("B43").Copy()
$ws.Range("B" + $r1).PasteSpecial(-4122)
$ws.Range("B" + $r2).PasteSpecial(-4122)
$ws.Range("B44").Copy()
$ws.Range("B43").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Re-apply the table's descending sort on "Puntuación total" (column C) so the
# two new rows land in their correct rank position, same as every other entry.
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("C3:C" + ($lo.Range.Row + $lo.Range.Rows.Count - 1)), 0, 2)
$lo.Sort.Header = 1
$lo.Sort.Apply()

$ws.Range("C93").Select()
